$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "16.03.2023"
$ws.Range("C3").Value = "https://gitlab.intra.infineon.com/digital-reference/order_management/-/commit/47cbb7faff327805f4f0ae6f71ccbeec1e086e96"
$ws.Range("D3").Value = "b67d5d8d24cae366f52cd197bc9be1e731229e148dc25959342bb7020c3a6bd0"
